$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Lgi2"
$ws.Cells.Item(2, 3).Value = "Adam23"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.1188046666666667
$ws.Cells.Item(2, 8).Value = 0.356414
$ws.Cells.Item(2, 9).Value = 0.00979853232878679
$ws.Cells.Item(2, 10).Value = 0.009798532328786792
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.03926266666666667
$ws.Cells.Item(2, 14).Value = 0.117788
$ws.Cells.Item(2, 15).Value = 0.005313231574131687
$ws.Cells.Item(2, 16).Value = 0.005313231574131686
$ws.Cells.Item(2, 17).Value = 0.004664588025777778
$ws.Cells.Item(2, 18).Value = 0.041981292232
$ws.Cells.Item(2, 19).Value = 0.00005206187134946006
$ws.Cells.Item(2, 20).Value = 0.00005206187134946006

# Row 3
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Lgi2"
$ws.Cells.Item(3, 3).Value = "Adam23"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.1188046666666667
$ws.Cells.Item(3, 8).Value = 0.356414
$ws.Cells.Item(3, 9).Value = 0.00979853232878679
$ws.Cells.Item(3, 10).Value = 0.009798532328786792
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 4.402094666666667
$ws.Cells.Item(3, 14).Value = 13.206284
$ws.Cells.Item(3, 15).Value = 0.5957147173375057
$ws.Cells.Item(3, 16).Value = 0.5957147173375056
$ws.Cells.Item(3, 17).Value = 0.5229893895084444
$ws.Cells.Item(3, 18).Value = 4.706904505576
$ws.Cells.Item(3, 19).Value = 0.005837129916565634
$ws.Cells.Item(3, 20).Value = 0.005837129916565634

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Lgi2"
$ws.Cells.Item(4, 3).Value = "Adam23"
$ws.Cells.Item(4, 4).Value = "ECs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.1188046666666667
$ws.Cells.Item(4, 8).Value = 0.356414
$ws.Cells.Item(4, 9).Value = 0.00979853232878679
$ws.Cells.Item(4, 10).Value = 0.009798532328786792
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 2.948244666666667
$ws.Cells.Item(4, 14).Value = 8.844734
$ws.Cells.Item(4, 15).Value = 0.3989720510883627
$ws.Cells.Item(4, 16).Value = 0.3989720510883626
$ws.Cells.Item(4, 17).Value = 0.3502652248751112
$ws.Cells.Item(4, 18).Value = 3.152387023876
$ws.Cells.Item(4, 19).Value = 0.003909340540871696
$ws.Cells.Item(4, 20).Value = 0.003909340540871696

# Row 5
$ws.Cells.Item(5, 1).Value = "sCs"
$ws.Cells.Item(5, 2).Value = "Lgi2"
$ws.Cells.Item(5, 3).Value = "Adam23"
$ws.Cells.Item(5, 4).Value = "FAPs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 11.502271
$ws.Cells.Item(5, 8).Value = 34.50681299999999
$ws.Cells.Item(5, 9).Value = 0.9486611714015168
$ws.Cells.Item(5, 10).Value = 0.948661171401517
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.03926266666666667
$ws.Cells.Item(5, 14).Value = 0.117788
$ws.Cells.Item(5, 15).Value = 0.005313231574131687
$ws.Cells.Item(5, 16).Value = 0.005313231574131686
$ws.Cells.Item(5, 17).Value = 0.4516098321826666
$ws.Cells.Item(5, 18).Value = 4.064488489643999
$ws.Cells.Item(5, 19).Value = 0.005040456489043292
$ws.Cells.Item(5, 20).Value = 0.005040456489043292

# Row 6
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Lgi2"
$ws.Cells.Item(6, 3).Value = "Adam23"
$ws.Cells.Item(6, 4).Value = "sCs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 11.502271
$ws.Cells.Item(6, 8).Value = 34.50681299999999
$ws.Cells.Item(6, 9).Value = 0.9486611714015168
$ws.Cells.Item(6, 10).Value = 0.948661171401517
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 4.402094666666667
$ws.Cells.Item(6, 14).Value = 13.206284
$ws.Cells.Item(6, 15).Value = 0.5957147173375057
$ws.Cells.Item(6, 16).Value = 0.5957147173375056
$ws.Cells.Item(6, 17).Value = 50.63408582365466
$ws.Cells.Item(6, 18).Value = 455.706772412892
$ws.Cells.Item(6, 19).Value = 0.5651314215705217
$ws.Cells.Item(6, 20).Value = 0.5651314215705217

# Row 7
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Lgi2"
$ws.Cells.Item(7, 3).Value = "Adam23"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 11.502271
$ws.Cells.Item(7, 8).Value = 34.50681299999999
$ws.Cells.Item(7, 9).Value = 0.9486611714015168
$ws.Cells.Item(7, 10).Value = 0.948661171401517
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.948244666666667
$ws.Cells.Item(7, 14).Value = 8.844734
$ws.Cells.Item(7, 15).Value = 0.3989720510883627
$ws.Cells.Item(7, 16).Value = 0.3989720510883626
$ws.Cells.Item(7, 17).Value = 33.91150913030467
$ws.Cells.Item(7, 18).Value = 305.203582172742
$ws.Cells.Item(7, 19).Value = 0.378489293341952
$ws.Cells.Item(7, 20).Value = 0.3784892933419519

# Row 8
$ws.Cells.Item(8, 1).Value = "ECs"
$ws.Cells.Item(8, 2).Value = "Lgi2"
$ws.Cells.Item(8, 3).Value = "Adam23"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.5036653333333333
$ws.Cells.Item(8, 8).Value = 1.510996
$ws.Cells.Item(8, 9).Value = 0.04154029626969626
$ws.Cells.Item(8, 10).Value = 0.04154029626969627
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.03926266666666667
$ws.Cells.Item(8, 14).Value = 0.117788
$ws.Cells.Item(8, 15).Value = 0.005313231574131687
$ws.Cells.Item(8, 16).Value = 0.005313231574131686
$ws.Cells.Item(8, 17).Value = 0.01977524409422222
$ws.Cells.Item(8, 18).Value = 0.177977196848
$ws.Cells.Item(8, 19).Value = 0.0002207132137389349
$ws.Cells.Item(8, 20).Value = 0.0002207132137389349

# Row 9
$ws.Cells.Item(9, 1).Value = "ECs"
$ws.Cells.Item(9, 2).Value = "Lgi2"
$ws.Cells.Item(9, 3).Value = "Adam23"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.5036653333333333
$ws.Cells.Item(9, 8).Value = 1.510996
$ws.Cells.Item(9, 9).Value = 0.04154029626969626
$ws.Cells.Item(9, 10).Value = 0.04154029626969627
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 4.402094666666667
$ws.Cells.Item(9, 14).Value = 13.206284
$ws.Cells.Item(9, 15).Value = 0.5957147173375057
$ws.Cells.Item(9, 16).Value = 0.5957147173375056
$ws.Cells.Item(9, 17).Value = 2.217182477651555
$ws.Cells.Item(9, 18).Value = 19.954642298864
$ws.Cells.Item(9, 19).Value = 0.02474616585041835
$ws.Cells.Item(9, 20).Value = 0.02474616585041835

# Row 10
$ws.Cells.Item(10, 1).Value = "ECs"
$ws.Cells.Item(10, 2).Value = "Lgi2"
$ws.Cells.Item(10, 3).Value = "Adam23"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.5036653333333333
$ws.Cells.Item(10, 8).Value = 1.510996
$ws.Cells.Item(10, 9).Value = 0.04154029626969626
$ws.Cells.Item(10, 10).Value = 0.04154029626969627
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 2.948244666666667
$ws.Cells.Item(10, 14).Value = 8.844734
$ws.Cells.Item(10, 15).Value = 0.3989720510883627
$ws.Cells.Item(10, 16).Value = 0.3989720510883626
$ws.Cells.Item(10, 17).Value = 1.484928632784889
$ws.Cells.Item(10, 18).Value = 13.364357695064
$ws.Cells.Item(10, 19).Value = 0.01657341720553898
$ws.Cells.Item(10, 20).Value = 0.01657341720553898

